$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly crypto price/volume refresh (GitHub Actions data pull).
# Force text storage for every touched cell (NumberFormat "@") so
# numeric-looking strings (e.g. '225.34') are not silently coerced
# into real numbers by the COM value-assignment type inference -
# then reset the style back to Normal so no stray style index is
# left behind on the cell (matches original formatting exactly).
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '34.080.91'
Set-TextValue 'E2' '  +0.09%  '
Set-TextValue 'D3' '1.779.93'
Set-TextValue 'E3' '  -0.44%  '
Set-TextValue 'E4' '  +0.22%  '
Set-TextValue 'D5' '225.34'
Set-TextValue 'E5' '  -0.61%  '
Set-TextValue 'E6' '  -0.08%  '
Set-TextValue 'E7' '  +0.22%  '
Set-TextValue 'D8' '31.79'
Set-TextValue 'E8' '  -1.27%  '
Set-TextValue 'E9' '  -1.52%  '
Set-TextValue 'E10' '  +0.30%  '
Set-TextValue 'D11' '0.0947'
Set-TextValue 'E11' '  +0.85%  '
Set-TextValue 'D12' '2.036.85'
Set-TextValue 'E12' '  -0.40%  '
Set-TextValue 'D13' '1.784.74'
Set-TextValue 'E13' '  -0.22%  '
Set-TextValue 'D14' '10.91'
Set-TextValue 'E14' '  -3.43%  '
Set-TextValue 'D15' '34.079.75'
Set-TextValue 'E15' '  +0.19%  '
Set-TextValue 'E16' '  +0.28%  '
Set-TextValue 'E17' '  -0.18%  '
Set-TextValue 'D18' '67.58'
Set-TextValue 'E18' '  -0.31%  '
Set-TextValue 'D19' '244.70'
Set-TextValue 'E19' '  +0.85%  '
Set-TextValue 'D20' '0.0₃0787'
Set-TextValue 'E20' '  +1.77%  '
Set-TextValue 'B21' 'Dai'
Set-TextValue 'C21' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue 'D21' '1.00'
Set-TextValue 'E21' '  +0.29%  '
Set-TextValue 'B22' 'Avalanche'
Set-TextValue 'C22' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue 'D22' '10.89'
Set-TextValue 'E22' '  +1.72%  '
Set-TextValue 'E23' '  +0.21%  '
Set-TextValue 'E24' '  -1.06%  '
Set-TextValue 'D25' '161.53'
Set-TextValue 'E25' '  -0.32%  '
Set-TextValue 'E26' '  -0.69%  '
Set-TextValue 'D27' '16.23'
Set-TextValue 'E27' '  +0.16%  '
Set-TextValue 'D28' '0.113'
Set-TextValue 'E28' '  +0.75%  '
Set-TextValue 'E29' '  +0.30%  '
Set-TextValue 'E30' '  -0.57%  '
Set-TextValue 'D31' '0.0516'
Set-TextValue 'E31' '  +0.10%  '
Set-TextValue 'E32' '  +1.57%  '
Set-TextValue 'E33' '  +2.56%  '
Set-TextValue 'E34' '  -2.16%  '
Set-TextValue 'D35' '1.445.04'
Set-TextValue 'E35' '  +3.46%  '
Set-TextValue 'D36' '2.45'
Set-TextValue 'E36' '  +3.70%  '
Set-TextValue 'D37' '0.651'
Set-TextValue 'E37' '  -0.26%  '
Set-TextValue 'E38' '  +1.17%  '
Set-TextValue 'E39' '  -0.42%  '
Set-TextValue 'E40' '  +1.33%  '
Set-TextValue 'D41' '80.27'
Set-TextValue 'E41' '  +0.20%  '
Set-TextValue 'E42' '  +1.47%  '
Set-TextValue 'D43' '0.915'
Set-TextValue 'E43' '  -0.53%  '
Set-TextValue 'E44' '  +0.02%  '
Set-TextValue 'D45' '0.0518'
Set-TextValue 'E45' '  +2.14%  '
Set-TextValue 'E46' '  -0.57%  '
Set-TextValue 'E47' '  -0.24%  '
Set-TextValue 'D48' '1.937.87'
Set-TextValue 'E48' '  -0.41%  '
Set-TextValue 'D49' '104.21'
Set-TextValue 'E49' '  -3.11%  '
Set-TextValue 'E50' '  +0.22%  '
Set-TextValue 'D51' '0.0₆0130'
Set-TextValue 'E51' '  -6.38%  '
